$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(22)
$p2 = $d.Paragraphs.Item(23)
$r = $d.Range($p1.Range.Start, $p2.Range.End)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr><w:t>Sprint 6 Backlog Items:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph" /><w:numPr><w:ilvl w:val="0" /><w:numId w:val="2" /></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr><w:t>Find and do a Demo / Presentation if you have not done so already</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr><w:t xml:space="preserve">The whole team hammered down and decided upon which tasks to tackle in order to fulfill these last few requirements of the course. Jace presented his topic to us and then to the course, I recorded the video and shared with my colleagues prior to submitting it, and a few others. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph" /><w:numPr><w:ilvl w:val="0" /><w:numId w:val="2" /></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr><w:t>Present Klump Project</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr><w:t xml:space="preserve">This corresponds to the previous backlog item, but was a very task considering the Klump Product was the big Kahuna of this semester. Julian Moses did a great job showcasing our hard work and even got </w:t></w:r><w:proofErr w:type="gramStart" /><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr><w:t>a the</w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr><w:t xml:space="preserve"> whole class laughing by telling a joke. Overall, a huge success.</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi" /><w:sz w:val="28" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center" /><w:rPr><w:rFonts w:ascii="Edwardian Script ITC" w:hAnsi="Edwardian Script ITC" w:cstheme="minorHAnsi" /><w:i /><w:sz w:val="28" /><w:u w:val="single" /></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Edwardian Script ITC" w:hAnsi="Edwardian Script ITC" w:cstheme="minorHAnsi" /><w:i /><w:sz w:val="28" /><w:u w:val="single" /></w:rPr><w:t>Michael Pedzimaz</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack" /><w:bookmarkEnd w:id="0" /></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center" /><w:rPr><w:rFonts w:ascii="Edwardian Script ITC" w:hAnsi="Edwardian Script ITC" w:cstheme="minorHAnsi" /><w:sz w:val="28" /><w:u w:val="single" /></w:rPr></w:pPr></w:p>'

$r.InsertXML($xml)
